$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text columns (G-K) need a leading quote-prefix so that
# Excel stores them as text (matching the existing t="str"-style data)
# rather than auto-converting them to numbers.

# Row 17
$ws.Cells.Item(17, 1).Value = " Abu Dhabi"
$ws.Cells.Item(17, 2).Value = " October 25 2020"
$ws.Cells.Item(17, 3).Value = "Royals won by 8 wickets (with 10 balls remaining)"
$ws.Cells.Item(17, 4).Value = "Mumbai Indians"
$ws.Cells.Item(17, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(17, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(17, 7).Value = "'40"
$ws.Cells.Item(17, 8).Value = "'26"
$ws.Cells.Item(17, 9).Value = "'4"
$ws.Cells.Item(17, 10).Value = "'1"
$ws.Cells.Item(17, 11).Value = "'153.84"

# Row 18
$ws.Cells.Item(18, 1).Value = " Abu Dhabi"
$ws.Cells.Item(18, 2).Value = " September 19 2020"
$ws.Cells.Item(18, 3).Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Cells.Item(18, 4).Value = "Mumbai Indians"
$ws.Cells.Item(18, 5).Value = "Chennai Super Kings"
$ws.Cells.Item(18, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(18, 7).Value = "'17"
$ws.Cells.Item(18, 8).Value = "'16"
$ws.Cells.Item(18, 9).Value = "'2"
$ws.Cells.Item(18, 10).Value = "'0"
$ws.Cells.Item(18, 11).Value = "'106.25"

# Row 19
$ws.Cells.Item(19, 1).Value = " Abu Dhabi"
$ws.Cells.Item(19, 2).Value = " October 28 2020"
$ws.Cells.Item(19, 3).Value = "Mumbai won by 5 wickets (with 5 balls remaining)"
$ws.Cells.Item(19, 4).Value = "Mumbai Indians"
$ws.Cells.Item(19, 5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(19, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(19, 7).Value = "'79"
$ws.Cells.Item(19, 8).Value = "'43"
$ws.Cells.Item(19, 9).Value = "'10"
$ws.Cells.Item(19, 10).Value = "'3"
$ws.Cells.Item(19, 11).Value = "'183.72"

# Row 20
$ws.Cells.Item(20, 1).Value = " Abu Dhabi"
$ws.Cells.Item(20, 2).Value = " September 23 2020"
$ws.Cells.Item(20, 3).Value = "Mumbai won by 49 runs"
$ws.Cells.Item(20, 4).Value = "Mumbai Indians"
$ws.Cells.Item(20, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(20, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(20, 7).Value = "'47"
$ws.Cells.Item(20, 8).Value = "'28"
$ws.Cells.Item(20, 9).Value = "'6"
$ws.Cells.Item(20, 10).Value = "'1"
$ws.Cells.Item(20, 11).Value = "'167.85"

# Row 21
$ws.Cells.Item(21, 1).Value = " Abu Dhabi"
$ws.Cells.Item(21, 2).Value = " October 11 2020"
$ws.Cells.Item(21, 3).Value = "Mumbai won by 5 wickets (with 2 balls remaining)"
$ws.Cells.Item(21, 4).Value = "Mumbai Indians"
$ws.Cells.Item(21, 5).Value = "Delhi Capitals"
$ws.Cells.Item(21, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(21, 7).Value = "'53"
$ws.Cells.Item(21, 8).Value = "'32"
$ws.Cells.Item(21, 9).Value = "'6"
$ws.Cells.Item(21, 10).Value = "'1"
$ws.Cells.Item(21, 11).Value = "'165.62"

# Row 22
$ws.Cells.Item(22, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(22, 2).Value = " November 05 2020"
$ws.Cells.Item(22, 3).Value = "Mumbai won by 57 runs"
$ws.Cells.Item(22, 4).Value = "Mumbai Indians"
$ws.Cells.Item(22, 5).Value = "Delhi Capitals"
$ws.Cells.Item(22, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(22, 7).Value = "'51"
$ws.Cells.Item(22, 8).Value = "'38"
$ws.Cells.Item(22, 9).Value = "'6"
$ws.Cells.Item(22, 10).Value = "'2"
$ws.Cells.Item(22, 11).Value = "'134.21"

# Row 23
$ws.Cells.Item(23, 1).Value = " Abu Dhabi"
$ws.Cells.Item(23, 2).Value = " October 01 2020"
$ws.Cells.Item(23, 3).Value = "Mumbai won by 48 runs"
$ws.Cells.Item(23, 4).Value = "Mumbai Indians"
$ws.Cells.Item(23, 5).Value = "Kings XI Punjab"
$ws.Cells.Item(23, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(23, 7).Value = "'10"
$ws.Cells.Item(23, 8).Value = "'7"
$ws.Cells.Item(23, 9).Value = "'2"
$ws.Cells.Item(23, 10).Value = "'0"
$ws.Cells.Item(23, 11).Value = "'142.85"

# Row 24
$ws.Cells.Item(24, 1).Value = " Sharjah"
$ws.Cells.Item(24, 2).Value = " October 04 2020"
$ws.Cells.Item(24, 3).Value = "Mumbai won by 34 runs"
$ws.Cells.Item(24, 4).Value = "Mumbai Indians"
$ws.Cells.Item(24, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(24, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(24, 7).Value = "'27"
$ws.Cells.Item(24, 8).Value = "'18"
$ws.Cells.Item(24, 9).Value = "'6"
$ws.Cells.Item(24, 10).Value = "'0"
$ws.Cells.Item(24, 11).Value = "'150.00"

# Row 25
$ws.Cells.Item(25, 1).Value = " Sharjah"
$ws.Cells.Item(25, 2).Value = " November 03 2020"
$ws.Cells.Item(25, 3).Value = "Sunrisers won by 10 wickets (with 17 balls remaining)"
$ws.Cells.Item(25, 4).Value = "Mumbai Indians"
$ws.Cells.Item(25, 5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(25, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(25, 7).Value = "'36"
$ws.Cells.Item(25, 8).Value = "'29"
$ws.Cells.Item(25, 9).Value = "'5"
$ws.Cells.Item(25, 10).Value = "'0"
$ws.Cells.Item(25, 11).Value = "'124.13"

# Row 26
$ws.Cells.Item(26, 1).Value = " Abu Dhabi"
$ws.Cells.Item(26, 2).Value = " October 16 2020"
$ws.Cells.Item(26, 3).Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Cells.Item(26, 4).Value = "Mumbai Indians"
$ws.Cells.Item(26, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(26, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(26, 7).Value = "'10"
$ws.Cells.Item(26, 8).Value = "'10"
$ws.Cells.Item(26, 9).Value = "'1"
$ws.Cells.Item(26, 10).Value = "'0"
$ws.Cells.Item(26, 11).Value = "'100.00"

# Row 27
$ws.Cells.Item(27, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(27, 2).Value = " October 31 2020"
$ws.Cells.Item(27, 3).Value = "Mumbai won by 9 wickets (with 34 balls remaining)"
$ws.Cells.Item(27, 4).Value = "Mumbai Indians"
$ws.Cells.Item(27, 5).Value = "Delhi Capitals"
$ws.Cells.Item(27, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(27, 7).Value = "'12"
$ws.Cells.Item(27, 8).Value = "'11"
$ws.Cells.Item(27, 9).Value = "'1"
$ws.Cells.Item(27, 10).Value = "'0"
$ws.Cells.Item(27, 11).Value = "'109.09"

# Row 28
$ws.Cells.Item(28, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(28, 2).Value = " November 10 2020"
$ws.Cells.Item(28, 3).Value = "Mumbai won by 5 wickets (with 8 balls remaining)"
$ws.Cells.Item(28, 4).Value = "Mumbai Indians"
$ws.Cells.Item(28, 5).Value = "Delhi Capitals"
$ws.Cells.Item(28, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(28, 7).Value = "'19"
$ws.Cells.Item(28, 8).Value = "'20"
$ws.Cells.Item(28, 9).Value = "'1"
$ws.Cells.Item(28, 10).Value = "'1"
$ws.Cells.Item(28, 11).Value = "'95.00"

# Row 29
$ws.Cells.Item(29, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(29, 2).Value = " September 28 2020"
$ws.Cells.Item(29, 3).Value = "Match tied (RCB won the one-over eliminator)"
$ws.Cells.Item(29, 4).Value = "Mumbai Indians"
$ws.Cells.Item(29, 5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(29, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(29, 7).Value = "'0"
$ws.Cells.Item(29, 8).Value = "'2"
$ws.Cells.Item(29, 9).Value = "'0"
$ws.Cells.Item(29, 10).Value = "'0"
$ws.Cells.Item(29, 11).Value = "'0.00"

# Row 30
$ws.Cells.Item(30, 1).Value = " Abu Dhabi"
$ws.Cells.Item(30, 2).Value = " October 06 2020"
$ws.Cells.Item(30, 3).Value = "Mumbai won by 57 runs"
$ws.Cells.Item(30, 4).Value = "Mumbai Indians"
$ws.Cells.Item(30, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(30, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(30, 7).Value = "'79"
$ws.Cells.Item(30, 8).Value = "'47"
$ws.Cells.Item(30, 9).Value = "'11"
$ws.Cells.Item(30, 10).Value = "'2"
$ws.Cells.Item(30, 11).Value = "'168.08"

# Row 31
$ws.Cells.Item(31, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(31, 2).Value = " October 18 2020"
$ws.Cells.Item(31, 3).Value = "Match tied (Kings XI won the one-over eliminator)"
$ws.Cells.Item(31, 4).Value = "Mumbai Indians"
$ws.Cells.Item(31, 5).Value = "Kings XI Punjab"
$ws.Cells.Item(31, 6).Value = "Suryakumar Yadav "
$ws.Cells.Item(31, 7).Value = "'0"
$ws.Cells.Item(31, 8).Value = "'4"
$ws.Cells.Item(31, 9).Value = "'0"
$ws.Cells.Item(31, 10).Value = "'0"
$ws.Cells.Item(31, 11).Value = "'0.00"

